$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Inhouse Data 2023" block (rows 12-20) had duplicate rows left over from
# pasting the same measurement multiple times. Row 13 is updated with the
# kon/koff/Kd figures that used to live on row 16, and the now-redundant
# rows 14-20 are cleaned out (14-17 become blank placeholder rows, 18-20 are
# removed outright).
$ws.Range("F13").Value = 1010000
$ws.Range("H13").Value = 0.00211
$ws.Range("J13").Value = 0.00000000529

# Wipe the leftover reference/method/ligand/receptor/parameter text in rows 14-20 ...
$ws.Range("A14:E20").ClearContents()
# ... and the leftover kon/koff/Kd numbers in rows 14-17 (these rows stay, just empty) ...
$ws.Range("F14:F17").ClearContents()
$ws.Range("H14:H17").ClearContents()
$ws.Range("J14:J17").ClearContents()

# ... while rows 18-20 are deleted completely, shifting nothing up into them.
$ws.Rows("18:20").Delete()

# Restore the view: zoomed to 73%, scrolled down a bit with C21 selected.
$win = $excel.ActiveWindow
$win.Zoom = 73
$null = $ws.Range("C21").Select()
